$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Remove the "Meta description" paragraph (2nd paragraph of the document).
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2) Insert a new paragraph ("Play Book of Kings Slot for Free - Review", bold)
#    right before the final "Prompt: ..." paragraph, i.e. after what is now the
#    second-to-last paragraph ("Lower frequency compared to other online slots").
$count = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs.Item($count - 1)
$insertPos = $precedingPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$newParaXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Kings Slot for Free - Review</w:t></w:r></w:p>'
[void]$insertRange.InsertXML($newParaXml)

# 3) Replace the text of the final "Prompt: ..." paragraph with the new
#    meta-description text, keeping its italic run formatting and leading
#    empty run intact.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$apost = [char]39
$newText = 'Read our Book of Kings slot review and play for free! Discover the game' + $apost + 's features, pros, and cons. An excellent option for experienced gamblers.'
$replacementXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p>'
[void]$lastPara.Range.InsertXML($replacementXml)
